# Login test and fix general tab with assert
# - Remove the unused trailing columns (Q:AW) from the "RS" sheet so the
#   sheet's used range shrinks back down to A1:P15.
# - Rename/retype a few header + data cells (ASCII-fy diacritics, swap a
#   handful of test-data values) and populate the new "Banka" column (N).
# - Leave the cursor parked on E3, matching the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RS")

# Drop the long tail of empty columns (Q:AW) that used to stretch the sheet
# out to AW. Deleting the whole columns shifts dimension back to A1:P15.
$ws.Range("Q1:AW1").EntireColumn.Delete()

# ---- Header row (row 1) ----
$ws.Range("H1").Value = "Drzava"
$ws.Range("J1").Value = "Dinamika placanja"
$ws.Range("M1").Value = "Metod placanja"

# ---- Row 2 ----
$ws.Range("G2").Value = "8"
$ws.Range("N2").Value = "0"

# ---- Row 3 ----
$ws.Range("J3").Value = "Polugodišnje"
$ws.Range("L3").Value = "Ne"
$ws.Range("M3").Value = "Trajni nalog"
$ws.Range("N3").Value = "BANCA INTESA AD BEOGRAD"

# ---- Row 4 ----
$ws.Range("E4").Value = "Riziko"
$ws.Range("J4").Value = "Mesečno"
$ws.Range("M4").Value = "Administrativna zabrana"
$ws.Range("N4").Value = "0"

# ---- Row 5 ----
$ws.Range("J5").Value = "Kvartalno"
$ws.Range("L5").Value = "Ne"
$ws.Range("M5").Value = "Nalog za uplatu premije"
$ws.Range("N5").Value = "0"

# ---- Row 6 ----
$ws.Range("N6").Value = "0"

# ---- Row 7 ----
$ws.Range("J7").Value = "Polugodišnje"
$ws.Range("M7").Value = "Trajni nalog"
$ws.Range("N7").Value = "BANCA INTESA AD BEOGRAD"

# ---- Row 8 ----
$ws.Range("G8").Value = "Doživotno"
$ws.Range("J8").Value = "Mesečno"
$ws.Range("M8").Value = "Administrativna zabrana"
$ws.Range("N8").Value = "0"

# ---- Row 9 ----
$ws.Range("J9").Value = "Kvartalno"
$ws.Range("L9").Value = "Ne"
$ws.Range("M9").Value = "Nalog za uplatu premije"
$ws.Range("N9").Value = "0"

# ---- Row 10 ----
$ws.Range("N10").Value = "0"

# ---- Row 11 ----
$ws.Range("E11").Value = "Riziko"
$ws.Range("J11").Value = "Polugodišnje"
$ws.Range("L11").Value = "Ne"
$ws.Range("M11").Value = "Trajni nalog"
$ws.Range("N11").Value = "BANCA INTESA AD BEOGRAD"

# ---- Row 12 ----
$ws.Range("J12").Value = "Mesečno"
$ws.Range("M12").Value = "Administrativna zabrana"
$ws.Range("N12").Value = "0"

# ---- Row 13 ----
$ws.Range("J13").Value = "Kvartalno"
$ws.Range("L13").Value = "Ne"
$ws.Range("M13").Value = "Nalog za uplatu premije"
$ws.Range("N13").Value = "0"

# ---- Row 14 ----
$ws.Range("N14").Value = "0"

# ---- Row 15 ----
$ws.Range("G15").Value = "Doživotno"
$ws.Range("J15").Value = "Polugodišnje"
$ws.Range("L15").Value = "Ne"
$ws.Range("M15").Value = "Trajni nalog"
$ws.Range("N15").Value = "BANCA INTESA AD BEOGRAD"

# Park the selection/active cell on E3, like the saved workbook.
$ws.Range("E3").Select() | Out-Null
